$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "It should not be displayed and should be as per parameter."

$ws.Range("F19").Value = $newText
$ws.Range("F20").Value = $newText
$ws.Range("F21").Value = $newText

$ws.Range("F46").Value = $newText
$ws.Range("F47").Value = $newText
$ws.Range("F48").Value = $newText

$ws.Range("F68").Value = $newText
$ws.Range("F69").Value = $newText
$ws.Range("F70").Value = $newText

$ws.Rows("68:70").AutoFit()

# Scroll the window so row 58 is the top visible row, then move the
# selection to match the saved view state (E77).
$excel.ActiveWindow.ScrollRow = 58
$ws.Range("E77").Select()

